$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4349165  # H17
$ws.Cells.Item(17, 10).Value = 4546808.5  # J17
$ws.Cells.Item(17, 12).Value = 13640425.5  # L17
$ws.Cells.Item(17, 14).Value = -13640761.5  # N17

$ws.Cells.Item(100, 8).Value = 1238.4348  # H100
$ws.Cells.Item(100, 9).Value = 975.25  # I100
$ws.Cells.Item(100, 11).Value = 975.25  # K100
$ws.Cells.Item(100, 13).Value = -434.25  # M100

$ws.Cells.Item(112, 8).Value = 3414.2856  # H112
$ws.Cells.Item(112, 10).Value = 3316.6667  # J112
$ws.Cells.Item(112, 12).Value = 9950.000100000001  # L112
$ws.Cells.Item(112, 14).Value = -12166.0001  # N112

$ws.Cells.Item(123, 8).Value = 80613.336  # H123
$ws.Cells.Item(123, 10).Value = 80613.336  # J123
$ws.Cells.Item(123, 12).Value = 80613.336  # L123
$ws.Cells.Item(123, 14).Value = -90413.336  # N123

$ws.Cells.Item(138, 8).Value = 1886.961  # H138
$ws.Cells.Item(138, 10).Value = 2355.8572  # J138
$ws.Cells.Item(138, 12).Value = 7067.571599999999  # L138
$ws.Cells.Item(138, 14).Value = -17347.5716  # N138

$ws.Cells.Item(141, 8).Value = 2436.25  # H141
$ws.Cells.Item(141, 9).Value = 2436.25  # I141
$ws.Cells.Item(141, 10).Value = 0  # J141
$ws.Cells.Item(141, 11).Value = 7308.75  # K141
$ws.Cells.Item(141, 12).Value = 0  # L141
$ws.Cells.Item(141, 13).Value = -2128.75  # M141
$ws.Cells.Item(141, 14).ClearContents()  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 50000  # H24
$ws.Cells.Item(24, 10).Value = 50000  # J24
$ws.Cells.Item(24, 12).Value = 50000  # L24
$ws.Cells.Item(24, 14).Value = -50748  # N24

$ws.Cells.Item(32, 8).Value = 3776.75  # H32
$ws.Cells.Item(32, 9).Value = 3506.1956  # I32
$ws.Cells.Item(32, 10).Value = 9999.5  # J32
$ws.Cells.Item(32, 11).Value = 3506.1956  # K32
$ws.Cells.Item(32, 12).Value = 9999.5  # L32
$ws.Cells.Item(32, 13).Value = -3219.1956  # M32
$ws.Cells.Item(32, 14).Value = -10573.5  # N32

$ws.Cells.Item(61, 8).Value = 6263  # H61
$ws.Cells.Item(61, 9).Value = 5307.6665  # I61
$ws.Cells.Item(61, 11).Value = 5307.6665  # K61
$ws.Cells.Item(61, 13).Value = -5095.6665  # M61

$ws.Cells.Item(100, 8).Value = 50000  # H100
$ws.Cells.Item(100, 10).Value = 50000  # J100
$ws.Cells.Item(100, 12).Value = 50000  # L100
$ws.Cells.Item(100, 14).Value = -52164  # N100

$ws.Cells.Item(102, 8).Value = 6929.9  # H102
$ws.Cells.Item(102, 9).Value = 4757  # I102
$ws.Cells.Item(102, 11).Value = 4757  # K102
$ws.Cells.Item(102, 13).Value = -3135  # M102

$ws.Cells.Item(110, 8).Value = 3750.4285  # H110
$ws.Cells.Item(110, 9).Value = 2192.3044  # I110
$ws.Cells.Item(110, 10).Value = 10917.8  # J110
$ws.Cells.Item(110, 11).Value = 2192.3044  # K110
$ws.Cells.Item(110, 12).Value = 10917.8  # L110
$ws.Cells.Item(110, 13).Value = -147.3044  # M110
$ws.Cells.Item(110, 14).Value = -15007.8  # N110

$ws.Cells.Item(132, 8).Value = 8208.25  # H132
$ws.Cells.Item(132, 9).Value = 7147.8887  # I132
$ws.Cells.Item(132, 10).Value = 11389.333  # J132
$ws.Cells.Item(132, 11).Value = 21443.6661  # K132
$ws.Cells.Item(132, 12).Value = 34167.999  # L132
$ws.Cells.Item(132, 13).Value = -18913.6661  # M132
$ws.Cells.Item(132, 14).Value = -39227.999  # N132

$ws.Cells.Item(136, 8).Value = 6263  # H136
$ws.Cells.Item(136, 9).Value = 5307.6665  # I136
$ws.Cells.Item(136, 11).Value = 15922.9995  # K136
$ws.Cells.Item(136, 13).Value = -13372.9995  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 3000  # H16
$ws.Cells.Item(16, 10).Value = 3000  # J16
$ws.Cells.Item(16, 12).Value = 3000  # L16
$ws.Cells.Item(16, 14).Value = -3340  # N16

$ws.Cells.Item(20, 8).Value = 3371.1  # H20
$ws.Cells.Item(20, 9).Value = 2536.4707  # I20
$ws.Cells.Item(20, 10).Value = 4462.5386  # J20
$ws.Cells.Item(20, 11).Value = 2536.4707  # K20
$ws.Cells.Item(20, 12).Value = 4462.5386  # L20
$ws.Cells.Item(20, 13).Value = -2289.4707  # M20
$ws.Cells.Item(20, 14).Value = -4956.5386  # N20

$ws.Cells.Item(99, 8).Value = 5549.684  # H99
$ws.Cells.Item(99, 9).Value = 4449.4546  # I99
$ws.Cells.Item(99, 10).Value = 7062.5  # J99
$ws.Cells.Item(99, 11).Value = 4449.4546  # K99
$ws.Cells.Item(99, 12).Value = 7062.5  # L99
$ws.Cells.Item(99, 13).Value = -2951.4546  # M99
$ws.Cells.Item(99, 14).Value = -10058.5  # N99

$ws.Cells.Item(105, 8).Value = 333.33334  # H105
$ws.Cells.Item(105, 9).Value = 333.33334  # I105
$ws.Cells.Item(105, 11).Value = 333.33334  # K105
$ws.Cells.Item(105, 13).Value = 1413.66666  # M105

$ws.Cells.Item(131, 8).Value = 70999.5  # H131
$ws.Cells.Item(131, 10).Value = 70999.5  # J131
$ws.Cells.Item(131, 12).Value = 70999.5  # L131
$ws.Cells.Item(131, 14).Value = -81079.5  # N131

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4996.1924  # H31
$ws.Cells.Item(31, 9).Value = 4482.3335  # I31
$ws.Cells.Item(31, 10).Value = 5436.643  # J31
$ws.Cells.Item(31, 11).Value = 4482.3335  # K31
$ws.Cells.Item(31, 12).Value = 5436.643  # L31
$ws.Cells.Item(31, 13).Value = -4187.3335  # M31
$ws.Cells.Item(31, 14).Value = -6026.643  # N31

$ws.Cells.Item(34, 8).Value = 4996.1924  # H34
$ws.Cells.Item(34, 9).Value = 4482.3335  # I34
$ws.Cells.Item(34, 10).Value = 5436.643  # J34
$ws.Cells.Item(34, 11).Value = 4482.3335  # K34
$ws.Cells.Item(34, 12).Value = 5436.643  # L34
$ws.Cells.Item(34, 13).Value = -4280.3335  # M34
$ws.Cells.Item(34, 14).Value = -5840.643  # N34

$ws.Cells.Item(58, 8).Value = 5275.514  # H58
$ws.Cells.Item(58, 9).Value = 3334.963  # I58
$ws.Cells.Item(58, 10).Value = 11824.875  # J58
$ws.Cells.Item(58, 11).Value = 3334.963  # K58
$ws.Cells.Item(58, 12).Value = 11824.875  # L58
$ws.Cells.Item(58, 13).Value = -3131.963  # M58
$ws.Cells.Item(58, 14).Value = -12230.875  # N58

$ws.Cells.Item(99, 8).Value = 3143.2144  # H99
$ws.Cells.Item(99, 9).Value = 2837.7273  # I99
$ws.Cells.Item(99, 10).Value = 4263.3335  # J99
$ws.Cells.Item(99, 11).Value = 2837.7273  # K99
$ws.Cells.Item(99, 12).Value = 4263.3335  # L99
$ws.Cells.Item(99, 13).Value = -1339.7273  # M99
$ws.Cells.Item(99, 14).Value = -7259.3335  # N99

$ws.Cells.Item(105, 8).Value = 2874.4  # H105
$ws.Cells.Item(105, 10).Value = 3349.3333  # J105
$ws.Cells.Item(105, 12).Value = 3349.3333  # L105
$ws.Cells.Item(105, 14).Value = -6843.3333  # N105

$ws.Cells.Item(126, 8).Value = 3143.2144  # H126
$ws.Cells.Item(126, 9).Value = 2837.7273  # I126
$ws.Cells.Item(126, 10).Value = 4263.3335  # J126
$ws.Cells.Item(126, 11).Value = 8513.1819  # K126
$ws.Cells.Item(126, 12).Value = 12790.0005  # L126
$ws.Cells.Item(126, 13).Value = -6043.1819  # M126
$ws.Cells.Item(126, 14).Value = -17730.0005  # N126

$ws.Cells.Item(132, 8).Value = 3453.2144  # H132
$ws.Cells.Item(132, 9).Value = 3064.182  # I132
$ws.Cells.Item(132, 11).Value = 9192.545999999998  # K132
$ws.Cells.Item(132, 13).Value = -6662.545999999998  # M132

$ws.Cells.Item(136, 8).Value = 5275.514  # H136
$ws.Cells.Item(136, 9).Value = 3334.963  # I136
$ws.Cells.Item(136, 10).Value = 11824.875  # J136
$ws.Cells.Item(136, 11).Value = 10004.889  # K136
$ws.Cells.Item(136, 12).Value = 35474.625  # L136
$ws.Cells.Item(136, 13).Value = -7454.889000000001  # M136
$ws.Cells.Item(136, 14).Value = -40574.625  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1568.5217  # H5
$ws.Cells.Item(5, 9).Value = 1519.8572  # I5
$ws.Cells.Item(5, 11).Value = 4559.571599999999  # K5
$ws.Cells.Item(5, 13).Value = -4447.571599999999  # M5

$ws.Cells.Item(36, 8).Value = 6750338  # H36
$ws.Cells.Item(36, 9).Value = 6750338  # I36
$ws.Cells.Item(36, 11).Value = 20251014  # K36
$ws.Cells.Item(36, 13).Value = -20250845  # M36

$ws.Cells.Item(96, 8).Value = 0  # H96
$ws.Cells.Item(96, 10).Value = 0  # J96
$ws.Cells.Item(96, 12).Value = 0  # L96
$ws.Cells.Item(96, 14).ClearContents()  # N96

$ws.Cells.Item(126, 8).Value = 12999.8  # H126
$ws.Cells.Item(126, 9).Value = 9999.5  # I126
$ws.Cells.Item(126, 11).Value = 29998.5  # K126
$ws.Cells.Item(126, 13).Value = -25058.5  # M126

$ws.Cells.Item(129, 8).Value = 22963260  # H129
$ws.Cells.Item(129, 9).Value = 41793148  # I129
$ws.Cells.Item(129, 10).Value = 1443385.8  # J129
$ws.Cells.Item(129, 11).Value = 125379444  # K129
$ws.Cells.Item(129, 12).Value = 4330157.4  # L129
$ws.Cells.Item(129, 13).Value = -125374444  # M129
$ws.Cells.Item(129, 14).Value = -4340157.4  # N129

$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 9).Value = 0  # I133
$ws.Cells.Item(133, 11).Value = 0  # K133
$ws.Cells.Item(133, 13).ClearContents()  # M133

$ws.Cells.Item(135, 8).Value = 1568.5217  # H135
$ws.Cells.Item(135, 9).Value = 1519.8572  # I135
$ws.Cells.Item(135, 11).Value = 13678.7148  # K135
$ws.Cells.Item(135, 13).Value = -11143.7148  # M135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15272.417  # H70
$ws.Cells.Item(70, 9).Value = 12249.5  # I70
$ws.Cells.Item(70, 11).Value = 12249.5  # K70
$ws.Cells.Item(70, 13).Value = -11979.5  # M70

$ws.Cells.Item(73, 8).Value = 15272.417  # H73
$ws.Cells.Item(73, 9).Value = 12249.5  # I73
$ws.Cells.Item(73, 11).Value = 12249.5  # K73
$ws.Cells.Item(73, 13).Value = -11313.5  # M73

$ws.Cells.Item(113, 8).Value = 365287.47  # H113
$ws.Cells.Item(113, 9).Value = 501421.75  # I113
$ws.Cells.Item(113, 10).Value = 2262.6667  # J113
$ws.Cells.Item(113, 11).Value = 501421.75  # K113
$ws.Cells.Item(113, 12).Value = 2262.6667  # L113
$ws.Cells.Item(113, 13).Value = -499251.75  # M113
$ws.Cells.Item(113, 14).Value = -6602.6667  # N113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3587.25  # H7
$ws.Cells.Item(7, 10).Value = 3666.3333  # J7
$ws.Cells.Item(7, 12).Value = 3666.3333  # L7
$ws.Cells.Item(7, 14).Value = -3890.3333  # N7

$ws.Cells.Item(22, 8).Value = 1397.375  # H22
$ws.Cells.Item(22, 9).Value = 1396.5  # I22
$ws.Cells.Item(22, 10).Value = 1400  # J22
$ws.Cells.Item(22, 11).Value = 1396.5  # K22
$ws.Cells.Item(22, 12).Value = 1400  # L22
$ws.Cells.Item(22, 13).Value = -1101.5  # M22
$ws.Cells.Item(22, 14).Value = -1990  # N22

$ws.Cells.Item(27, 8).Value = 1397.375  # H27
$ws.Cells.Item(27, 9).Value = 1396.5  # I27
$ws.Cells.Item(27, 10).Value = 1400  # J27
$ws.Cells.Item(27, 11).Value = 1396.5  # K27
$ws.Cells.Item(27, 12).Value = 1400  # L27
$ws.Cells.Item(27, 13).Value = -1289.5  # M27
$ws.Cells.Item(27, 14).Value = -1614  # N27

$ws.Cells.Item(54, 8).Value = 0  # H54
$ws.Cells.Item(54, 10).Value = 0  # J54
$ws.Cells.Item(54, 12).Value = 0  # L54
$ws.Cells.Item(54, 14).ClearContents()  # N54

$ws.Cells.Item(98, 8).Value = 0  # H98
$ws.Cells.Item(98, 10).Value = 0  # J98
$ws.Cells.Item(98, 12).Value = 0  # L98
$ws.Cells.Item(98, 14).ClearContents()  # N98

$ws.Cells.Item(126, 8).Value = 3587.25  # H126
$ws.Cells.Item(126, 10).Value = 3666.3333  # J126
$ws.Cells.Item(126, 12).Value = 10998.9999  # L126
$ws.Cells.Item(126, 14).Value = -15938.9999  # N126

$ws.Cells.Item(132, 8).Value = 8760.303  # H132
$ws.Cells.Item(132, 9).Value = 9054.058000000001  # I132
$ws.Cells.Item(132, 11).Value = 27162.174  # K132
$ws.Cells.Item(132, 13).Value = -24632.174  # M132

$ws.Cells.Item(136, 8).Value = 3928.6765  # H136
$ws.Cells.Item(136, 9).Value = 3601.138  # I136
$ws.Cells.Item(136, 10).Value = 5828.4  # J136
$ws.Cells.Item(136, 11).Value = 10803.414  # K136
$ws.Cells.Item(136, 12).Value = 17485.2  # L136
$ws.Cells.Item(136, 13).Value = -8253.414000000001  # M136
$ws.Cells.Item(136, 14).Value = -22585.2  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 1424437.6  # H5
$ws.Cells.Item(5, 10).Value = 1424437.6  # J5
$ws.Cells.Item(5, 12).Value = 1424437.6  # L5
$ws.Cells.Item(5, 14).Value = -1424661.6  # N5

$ws.Cells.Item(81, 8).Value = 2129.5  # H81
$ws.Cells.Item(81, 9).Value = 899  # I81
$ws.Cells.Item(81, 10).Value = 8282  # J81
$ws.Cells.Item(81, 11).Value = 1798  # K81
$ws.Cells.Item(81, 12).Value = 16564  # L81
$ws.Cells.Item(81, 13).Value = -737  # M81
$ws.Cells.Item(81, 14).Value = -18686  # N81

$ws.Cells.Item(84, 8).Value = 2129.5  # H84
$ws.Cells.Item(84, 9).Value = 899  # I84
$ws.Cells.Item(84, 10).Value = 8282  # J84
$ws.Cells.Item(84, 11).Value = 8990  # K84
$ws.Cells.Item(84, 12).Value = 82820  # L84
$ws.Cells.Item(84, 13).Value = -3686  # M84
$ws.Cells.Item(84, 14).Value = -93428  # N84
